$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3d_object_library")

# Insert 3 new rows above row 11 (rows shift down; old row 11 becomes row 14, etc.)
$ws.Rows("11:13").Insert()

# Match the "-1" header-row style (quote-prefixed integer format) used by rows 3-10
$ws.Range("A11:A13").Style = $ws.Range("A10").Style

# Row 11: Help 01 texture entry
$ws.Cells.Item(11, 1).Value = -1
$ws.Cells.Item(11, 2).Value = "Help 01"
$ws.Cells.Item(11, 3).Value = "texture"
$ws.Cells.Item(11, 4).Value = "help_01.png"
for ($c = 5; $c -le 32; $c++) {
    $ws.Cells.Item(11, $c).Value = "*"
}
$ws.Cells.Item(11, 33).Value = "0"

# Row 12: Help 02 texture entry
$ws.Cells.Item(12, 1).Value = -1
$ws.Cells.Item(12, 2).Value = "Help 02"
$ws.Cells.Item(12, 3).Value = "texture"
$ws.Cells.Item(12, 4).Value = "help_02.png"
for ($c = 5; $c -le 32; $c++) {
    $ws.Cells.Item(12, $c).Value = "*"
}
$ws.Cells.Item(12, 33).Value = "0"

# Row 13: Help 03 texture entry
$ws.Cells.Item(13, 1).Value = -1
$ws.Cells.Item(13, 2).Value = "Help 03"
$ws.Cells.Item(13, 3).Value = "texture"
$ws.Cells.Item(13, 4).Value = "help_03.png"
for ($c = 5; $c -le 32; $c++) {
    $ws.Cells.Item(13, $c).Value = "*"
}
$ws.Cells.Item(13, 33).Value = "0"

$ws.Range("D16").Select()
